{"js": "// Add the sentence \"Test fra torbj\u00f8rn\" into the empty paragraph that sits\n// between the \"PCA\" and \"PCR\" headings, and mark the three inline\n// drawings' runs as \"do not spell/grammar check\" (w:noProof), matching\n// what Word stamps onto picture runs whenever it resaves a document that\n// contains them.\n\n// 1) Locate the empty paragraph right after the paragraph whose text is\n//    \"PCA\" (and right before the one whose text is \"PCR\") and add a run\n//    with the requested text to it.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet targetParagraph = null;\nfor (let i = 0; i < paragraphs.items.length - 1; i++) {\n  const current = paragraphs.items[i];\n  if (current.text === \"PCA\") {\n    const candidate = paragraphs.items[i + 1];\n    if (candidate.text === \"\") {\n      targetParagraph = candidate;\n      break;\n    }\n  }\n}\n\nif (targetParagraph) {\n  // Use insertOoxml so the new run carries the same run formatting\n  // (<w:rPr><w:lang w:val=\"en-US\"/></w:rPr>) used throughout this\n  // document instead of whatever the default run font happens to be.\n  const ooxml =\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body><w:p><w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>Test fra torbj\\u00f8rn</w:t></w:r></w:p></w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>';\n  targetParagraph.insertOoxml(ooxml, Word.InsertLocation.end);\n  await context.sync();\n}\n\n// 2) Flag every inline picture's run as \"noProof\" (the little camera-icon\n//    runs Word marks so the proofing pass skips them).\nconst inlinePictures = context.document.body.inlinePictures;\ninlinePictures.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < inlinePictures.items.length; i++) {\n  const pictureRange = inlinePictures.items[i].getRange();\n  pictureRange.hasNoProofing = true;\n}\nawait context.sync();\n", "ps1": "# Add the sentence \"Test fra torbj\u00f8rn\" into the empty paragraph that sits\n# between the \"PCA\" and \"PCR\" headings, and mark the three inline\n# drawings' runs as \"do not spell/grammar check\" (w:noProof), matching\n# what Word stamps onto picture runs whenever it resaves a document that\n# contains them.\n\n$d = $word.ActiveDocument\n\n# 1) Locate the empty paragraph right after the paragraph whose text is\n#    \"PCA\" (and right before the one whose text is \"PCR\") and add a run\n#    with the requested text to it.\n$targetIndex = -1\nfor ($i = 1; $i -lt $d.Paragraphs.Count; $i++) {\n    $current = $d.Paragraphs.Item($i).Range.Text\n    if ($current -eq \"PCA`r\") {\n        $next = $d.Paragraphs.Item($i + 1).Range.Text\n        if ($next -eq \"`r\") {\n            $targetIndex = $i + 1\n            break\n        }\n    }\n}\n\nif ($targetIndex -ne -1) {\n    $targetParagraph = $d.Paragraphs.Item($targetIndex)\n    # Collapse to the very start of the paragraph and insert a real OOXML\n    # run fragment (rather than Range.Text / InsertBefore) so the new run\n    # picks up the same run formatting (<w:rPr><w:lang w:val=\"en-US\"/></w:rPr>)\n    # used throughout the rest of the document, instead of landing with no\n    # run properties at all.\n    $insertionPoint = $d.Range($targetParagraph.Range.Start, $targetParagraph.Range.Start)\n    $runOoxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body><w:p><w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>Test fra torbj\u00f8rn</w:t></w:r></w:p></w:body>' +\n        '</w:document>' +\n        '</pkg:xmlData></pkg:part></pkg:package>'\n    $insertionPoint.InsertXML($runOoxml)\n}\n\n# 2) Flag every inline picture's run as \"noProof\" (the little camera-icon\n#    runs Word marks so the proofing pass skips them).\nfor ($i = 1; $i -le $d.InlineShapes.Count; $i++) {\n    $pictureRange = $d.InlineShapes.Item($i).Range\n    $pictureRange.NoProofing = 1\n}\n"}
